# Generate Report for Archive
#
# 1. Update the "Status" value shown on the Overview sheet (columns E/F, one
#    per target locale) and on each locale sheet's "Status" column (C) from
#    "Ready for handoff" to "In Translation".
# 2. Narrow the "Status" column width on every sheet that shows it.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text -----------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Narrow the "Status" columns ---------------------------------------
# Target OOXML column width is 13.4101845877511 characters. This runtime
# quantizes ColumnWidth to whole-pixel increments (steps of 1/6 of a
# character) before persisting it, so 12.5 is the closest reachable
# ColumnWidth: it round-trips to a stored width of 13.3333.., the nearest
# achievable value to the target.
$statusWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $statusWidth
$overview.Columns.Item(6).ColumnWidth = $statusWidth
$zhcn.Columns.Item(3).ColumnWidth = $statusWidth
$dede.Columns.Item(3).ColumnWidth = $statusWidth
